$d = $word.ActiveDocument

function Get-ParaByText($needle) {
    $r = $d.Content.Duplicate
    $r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $r.Paragraphs(1)
}

function Set-ParaXml($para, $innerXml) {
    $wrapped = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($wrapped)
}

# 1. Title paragraph: drop the proofErr spellStart/spellEnd markers that bracketed "Kough".
$titlePara = Get-ParaByText("Kough")
$titleXml = '<w:p><w:pPr><w:pStyle w:val="Title"/><w:jc w:val="center"/><w:rPr><w:sz w:val="52"/><w:szCs w:val="52"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="52"/><w:szCs w:val="52"/><w:u w:val="single"/></w:rPr><w:t>Kough</w:t></w:r><w:r><w:rPr><w:sz w:val="52"/><w:szCs w:val="52"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> Privacy</w:t></w:r><w:r><w:rPr><w:sz w:val="52"/><w:szCs w:val="52"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> Policy</w:t></w:r></w:p>'
Set-ParaXml $titlePara $titleXml

# 2. Date paragraph: split "23/05/2022" into "2" + "5" + "/05/2022" runs.
$datePara = Get-ParaByText("23/05/2022")
$dateXml = '<w:p><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Our Privacy Policy was last updated on: </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>5</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>/05/2022</w:t></w:r></w:p>'
Set-ParaXml $datePara $dateXml

# 3. "To determine..." paragraph: drop the leading run that contains only a single space.
$toDeterminePara = Get-ParaByText("To determine where it is appropriate")
$toDetermineXml = '<w:p><w:r><w:t>To determine where it is appropriate a policy will be put in place, this policy will look at relevancy of the data and length of time it has been stored. Any data stored will also be stored securely to ensure we are the only people who can access it. Any breach will be reported to the correct body.</w:t></w:r></w:p>'
Set-ParaXml $toDeterminePara $toDetermineXml

# 4. The empty trailing "ListParagraph" bullet becomes a normal paragraph with new body text,
#    and a new paragraph is added straight after it (before the final blank paragraph).
$listPara = Get-ParaByText("Consensual storage of data to help and train")
$emptyBullet = $listPara.Next()
$bulletXml = '<w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>If</w:t></w:r><w:r><w:t xml:space="preserve">, any possible crimes </w:t></w:r><w:r><w:t>are</w:t></w:r><w:r><w:t xml:space="preserve"> intercepted by our team </w:t></w:r><w:r><w:t xml:space="preserve">during the time of recording </w:t></w:r><w:r><w:t xml:space="preserve">for analytic purposes </w:t></w:r><w:r><w:t>that can cause harm to an individual as per Investigatory Powers Act (2016) we’re obliged to report it to appropriate authorities.</w:t></w:r></w:p>'
Set-ParaXml $emptyBullet $bulletXml

$updatedBullet = Get-ParaByText("that can cause harm to an individual")
$updatedBullet.Range.InsertParagraphAfter()
$newPara = $updatedBullet.Next()
$newParaXml = '<w:p><w:r><w:t>The private policy operates on foundations of General Data Protection Act (2018).</w:t></w:r></w:p>'
Set-ParaXml $newPara $newParaXml
